$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for the error-table rows (B2:G11), reflecting the naive
# component forecaster bug fix: each prior row's B:F values shift down
# by one row (older observation), G (the count column) increments by
# one accordingly, and a new leading observation is computed for row 2.

$data = @(
    @{ Row = 2;  B = -0.03451060989511129; C = 1.472147746012473;  D = 8.907389105311308;  E = 2.984524937961033;  F = 3.051397271829799;  G = 23 }
    @{ Row = 3;  B = -0.1106822175073481;  C = 1.485348127483025;  D = 8.723650662140232;  E = 2.953582682462137;  F = 3.020964860020825;  G = 22 }
    @{ Row = 4;  B = -0.5713011636365682;  C = 1.053179700887051;  D = 4.026630981043165;  E = 2.006646700603563;  F = 1.971105846864699;  G = 21 }
    @{ Row = 5;  B = -0.166036682263726;   C = 0.6285186396184349; D = 0.8305896743800926; E = 0.911366926314584;  F = 0.9193942482462926; G = 20 }
    @{ Row = 6;  B = -0.1280209040135819;  C = 0.7182792096296192; D = 0.930922228823257;  E = 0.9648431109891685; F = 0.9825173072328391; G = 19 }
    @{ Row = 7;  B = -0.1613809308236712;  C = 0.6969901999678998; D = 0.7440183132057626; E = 0.8625649617308615; F = 0.8718992504770149; G = 18 }
    @{ Row = 8;  B = -0.06310423029847312; C = 0.6271081964451903; D = 0.5990212621699853; E = 0.7739646388369338; F = 0.7951283312657454; G = 17 }
    @{ Row = 9;  B = 0.003107267299192557; C = 0.5393758537464697; D = 0.4187293164217235; E = 0.6470929735530463; F = 0.6683070442536863; G = 16 }
    @{ Row = 10; B = 0.01079506215789655;  C = 0.4821603736089785; D = 0.3882225284189086; E = 0.6230750584150425; F = 0.6448471538599968; G = 15 }
    @{ Row = 11; B = 0.03588466164625848;  C = 0.5517454314036715; D = 0.4354142101678629; E = 0.6598592351159926; F = 0.6837549615721854; G = 14 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
    $ws.Cells.Item($r, 6).Value = $entry.F
    $ws.Cells.Item($r, 7).Value = $entry.G
}
